$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new master-location rows (postal code / Arabic rows for BNMR)
$ws.Range("A120").Value = 10113
$ws.Range("B120").Value = 10113
$ws.Range("C120").Value = 5
$ws.Range("D120").Value = "الرمز البريدي"
$ws.Range("E120").Value = "BNMR"
$ws.Range("F120").Value = "ara"
$ws.Range("G120").Value = $true
$ws.Range("H120").Value = "superadmin"
$ws.Range("I120").Value = "now()"

$ws.Range("A121").Value = 10114
$ws.Range("B121").Value = 10114
$ws.Range("C121").Value = 5
$ws.Range("D121").Value = "الرمز البريدي"
$ws.Range("E121").Value = "BNMR"
$ws.Range("F121").Value = "ara"
$ws.Range("G121").Value = $true
$ws.Range("H121").Value = "superadmin"
$ws.Range("I121").Value = "now()"

# Move/extend the selection to the full rows below the new data (A122:XFD1048576)
$excel.Goto($ws.Range("A122:XFD1048576"))
